# Append the new run-log row (row 69) to Sheet1, mirroring the existing
# rows' layout: A=Run UTC, B=Run IST, C=Status, D=Message, E=Chosen URL,
# F=Saved PDF, G=Rows Appended, H=Total Rows After.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 69

$ws.Cells.Item($row, 1).Value = "2025-08-28 13:03:36 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-28 18:33:36 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""

# Match the formatting used by the rest of the data rows (same cell
# style as the row above) by copying formats only.
$ws.Range("A68:H68").Copy()
$ws.Range("A$row`:H$row").PasteSpecial(-4122)
